$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update "Periodo Mora" (column E) for rows 16-22 so the periods run in
# chronological (ascending) order instead of the previous descending order.
$ws.Range("E16").Value = "1911"
$ws.Range("E17").Value = "1912"
$ws.Range("E18").Value = "2001"
$ws.Range("E19").Value = "2002"
$ws.Range("E20").Value = "2003"
$ws.Range("E21").Value = "2004"
$ws.Range("E22").Value = "2005"

# Update "Valor Mora" (column G) for rows 16-22 with the new amount.
$ws.Range("G16:G22").Value = 1423500
